# Weekly symbol-list refresh: updates the Price (D) and Volume(1h) (E)
# columns on Sheet1 with freshly scraped coinranking.com figures.
#
# All values in these columns are stored as plain text (inlineStr) rather
# than numbers/percentages, so each assignment is prefixed with a literal
# apostrophe (Excel's text-entry marker) to stop values such as "329.06"
# or "1.25%" from being auto-converted into a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.06"
$ws.Range("E2").Value = "'1.25%"
$ws.Range("D3").Value = "'41.19"
$ws.Range("D4").Value = "'5.617"
$ws.Range("E4").Value = "'-1.33%"
$ws.Range("D5").Value = "'0.08205"
$ws.Range("E5").Value = "'2.53%"
$ws.Range("D6").Value = "'8.746"
$ws.Range("E6").Value = "'1.38%"
$ws.Range("D7").Value = "'2.009"
$ws.Range("E7").Value = "'-0.11%"
$ws.Range("D8").Value = "'4.490"
$ws.Range("E8").Value = "'-0.21%"
$ws.Range("D10").Value = "'0.9197"
$ws.Range("E10").Value = "'-0.55%"
$ws.Range("D11").Value = "'0.1283"
$ws.Range("E11").Value = "'3.63%"
$ws.Range("D12").Value = "'0.1952"
$ws.Range("E12").Value = "'-1.25%"
$ws.Range("D13").Value = "'0.09314"
$ws.Range("E13").Value = "'0.90%"
$ws.Range("D14").Value = "'0.03896"
$ws.Range("E14").Value = "'7.45%"
$ws.Range("D15").Value = "'0.1061"
$ws.Range("E15").Value = "'0.98%"
$ws.Range("D16").Value = "'0.001308"
$ws.Range("E16").Value = "'1.06%"
$ws.Range("D17").Value = "'0.006102"
$ws.Range("E17").Value = "'-0.56%"
$ws.Range("E19").Value = "'2.82%"
$ws.Range("D21").Value = "'8.233"
$ws.Range("E21").Value = "'-5.36%"
$ws.Range("D22").Value = "'0.1366"
$ws.Range("E22").Value = "'-0.45%"
$ws.Range("D24").Value = "'0.04406"
$ws.Range("E24").Value = "'0.21%"
$ws.Range("D25").Value = "'0.001256"
$ws.Range("E25").Value = "'-0.43%"
$ws.Range("D26").Value = "'0.004312"
$ws.Range("E26").Value = "'-7.17%"
$ws.Range("E27").Value = "'4.33%"
$ws.Range("D39").Value = "'0.02775"
$ws.Range("E39").Value = "'11.46%"
$ws.Range("E40").Value = "'1.21%"
$ws.Range("D41").Value = "'0.007797"
$ws.Range("E41").Value = "'4.53%"
$ws.Range("E42").Value = "'0.70%"
$ws.Range("D43").Value = "'0.008948"
$ws.Range("E43").Value = "'-7.50%"
$ws.Range("E44").Value = "'2.55%"
$ws.Range("D45").Value = "'0.01219"
$ws.Range("E45").Value = "'9.39%"
$ws.Range("D46").Value = "'0.00006763"
$ws.Range("E46").Value = "'0.18%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("D48").Value = "'0.003191"
$ws.Range("E48").Value = "'7.37%"
$ws.Range("E49").Value = "'-0.47%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.03%"

# The apostrophe-prefixed entry leaves a "quote prefix" number format on
# each touched cell; clear that back off (per cell - ClearFormats on a
# multi-area union only affects the first area) so cells keep their
# original, unstyled formatting.
foreach ($ref in @("D2", "E2", "D3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E19", "D21", "E21", "D22", "E22", "D24", "E24", "D25", "E25", "D26", "E26", "E27", "D39", "E39", "E40", "D41", "E41", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "D51", "E51")) {
    $ws.Range($ref).ClearFormats()
}
